$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - copy formatting from existing header cell H1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns (rows 2-22)
$values = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(6, 7)
    5  = @(8, 9)
    6  = @(9, 9)
    7  = @(7, 8)
    8  = @(8, 9)
    9  = @(8, 9)
    10 = @(5, 7)
    11 = @(4, 7)
    12 = @(6, 8)
    13 = @(7, 9)
    14 = @(9, 9)
    15 = @(7, 8)
    16 = @(13, 13)
    17 = @(8, 9)
    18 = @(6, 6)
    19 = @(8, 9)
    20 = @(7, 9)
    21 = @(4, 4)
    22 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
